$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "standard_user"
$ws.Range("B2").Value = "secret_sauce"
$ws.Range("A3").Value = "problem_user"
$ws.Range("B3").Value = "secret_sauce"
$ws.Range("A4").Value = "performance_glitch_user"
$ws.Range("B4").Value = "secret_sauce"

$ws.Range("A5").Select()
